$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D (quarterly data shifts right by 2 columns)
$ws.Columns("D:E").Insert()

# Carry forward number formatting into the newly inserted columns from column F
# (the date header rows keep the date format, the data rows keep the numeric format)
$ws.Range("F7").Copy()
$ws.Range("D7:E7").PasteSpecial(-4122)
$ws.Range("F38").Copy()
$ws.Range("D38:E38").PasteSpecial(-4122)
$ws.Range("F80").Copy()
$ws.Range("D80:E80").PasteSpecial(-4122)

$ws.Range("F8:F35").Copy()
$ws.Range("D8:E35").PasteSpecial(-4122)
$ws.Range("F39:F102").Copy()
$ws.Range("D39:E102").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Populate the two new quarters of data (periods ending 2018-09-30 and 2018-12-31)
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 10100
$ws.Range("E8").Value = 9700
$ws.Range("D9").Value = 5600
$ws.Range("E9").Value = 4600
$ws.Range("D10").Value = 4500
$ws.Range("E10").Value = 5100
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = "NA"
$ws.Range("E14").Value = "NA"
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("D17").Value = 8300
$ws.Range("E17").Value = 5400
$ws.Range("D18").Value = 1800
$ws.Range("E18").Value = 4300
$ws.Range("D20").Value = 0
$ws.Range("E20").Value = 0
$ws.Range("D21").Value = "NA"
$ws.Range("E21").Value = "NA"
$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 0
$ws.Range("D23").Value = 1900
$ws.Range("E23").Value = 4400
$ws.Range("D24").Value = 1500
$ws.Range("E24").Value = "NA"
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 400
$ws.Range("E26").Value = 4400
$ws.Range("D27").Value = -500
$ws.Range("E27").Value = 3500
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("E29").Value = 0
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = 0
$ws.Range("E32").Value = 0
$ws.Range("D33").Value = -500
$ws.Range("E33").Value = 3500
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = -500
$ws.Range("E35").Value = 3500
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 7900
$ws.Range("E41").Value = 9900
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("D43").Value = 59800
$ws.Range("E43").Value = 18800
$ws.Range("D44").Value = 0
$ws.Range("E44").Value = 0
$ws.Range("D45").Value = 0
$ws.Range("E45").Value = 0
$ws.Range("D46").Value = 0
$ws.Range("E46").Value = 0
$ws.Range("D47").Value = 559200
$ws.Range("E47").Value = 573600
$ws.Range("D48").Value = 0
$ws.Range("E48").Value = 0
$ws.Range("D49").Value = 0
$ws.Range("E49").Value = 0
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 51500
$ws.Range("E52").Value = 77800
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 679400
$ws.Range("E54").Value = 680800
$ws.Range("D57").Value = 2100
$ws.Range("E57").Value = 100
$ws.Range("D58").Value = 0
$ws.Range("E58").Value = 0
$ws.Range("D59").Value = 3900
$ws.Range("E59").Value = 5300
$ws.Range("D60").Value = 0
$ws.Range("E60").Value = 0
$ws.Range("D61").Value = 19200
$ws.Range("E61").Value = 19500
$ws.Range("D62").Value = 0
$ws.Range("E62").Value = 0
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 529200
$ws.Range("E66").Value = 528700
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 37200
$ws.Range("E70").Value = 37200
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = -118600
$ws.Range("E72").Value = -116600
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 112900
$ws.Range("E76").Value = 114900
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = -500
$ws.Range("E81").Value = 3500
$ws.Range("D83").Value = 0
$ws.Range("E83").Value = 0
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 1200
$ws.Range("E89").Value = 5300
$ws.Range("D91").Value = 0
$ws.Range("E91").Value = 0
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -27400
$ws.Range("E94").Value = -213400
$ws.Range("D96").Value = -1400
$ws.Range("E96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = -2200
$ws.Range("E100").Value = 214500
$ws.Range("D101").Value = 0
$ws.Range("E101").Value = 0
$ws.Range("D102").Value = -28400
$ws.Range("E102").Value = 6300

# A few prior-period figures were also restated in this update
$ws.Range("F24").Value = "NA"
$ws.Range("G24").Value = "NA"
$ws.Range("H24").Value = "NA"
$ws.Range("I24").Value = "NA"
$ws.Range("J24").Value = "NA"

$ws.Range("F72").Value = -118700
$ws.Range("G72").Value = -94900
$ws.Range("H72").Value = -103000
$ws.Range("I72").Value = -107700
$ws.Range("J72").Value = -99200
